# Reclassify NBC AABM rows (Area 3/4/5/103/104/105 + their combined rows)
# to the new "NBC ISBM S ..." category, introducing three new shared
# strings: "NBC ISBM S FALL", "NBC ISBM S SPRING", "NBC ISBM S SUMMER".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Area")

$fallRows   = @(69, 70, 75, 76, 77, 80, 81, 82, 221)
$springRows = @(87, 88, 93, 94, 95, 98, 99, 100, 222)
$summerRows = @(105, 106, 111, 112, 113, 116, 117, 118, 223)

foreach ($r in $fallRows) {
    $ws.Cells.Item($r, 5).Value = "NBC ISBM S FALL"
}
foreach ($r in $springRows) {
    $ws.Cells.Item($r, 5).Value = "NBC ISBM S SPRING"
}
foreach ($r in $summerRows) {
    $ws.Cells.Item($r, 5).Value = "NBC ISBM S SUMMER"
}

# Move / restore the active selection to F95 (cosmetic UI state captured
# by the author's last save before committing).
$ws.Range("F95").Select()
